$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lop_hoc")

# The sheet lists classes ("lop_hoc"); two rows are being removed:
#   - the class dated 10/10/2020 (row 7)
#   - the class dated 2/10/2020 (originally row 9, becomes row 8 once the
#     first row above it is removed)
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(8).Delete()

# Deleting rows flattens the shared formulas that spanned the removed rows
# into ordinary per-cell formulas; re-apply them across their (now shorter)
# ranges so they collapse back into shared formulas like Excel normally
# keeps them.
$ws.Range("E3:E17").Formula = "=(D3+70)"
$ws.Range("G3:G17").Formula = "=TODAY()-D3"
$ws.Range("N4:N17").Formula = "=QUOTIENT(G4,7)"

# Match the post-delete selection Excel leaves behind (the whole row that
# now occupies row 7 is selected).
$ws.Activate()
$ws.Range("A7:XFD7").Select()
